$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.825.26"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.384.08"
$ws.Range("E3").Value = "  -2.45%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'543.33"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'140.69"
$ws.Range("E6").Value = "  -3.27%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.577"
$ws.Range("E8").Value = "  -5.05%  "
$ws.Range("D9").Value = "2.383.37"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("D12").Value = "'5.36"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "'0.343"
$ws.Range("E13").Value = "  -2.28%  "
$ws.Range("D14").Value = "'25.43"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").Value = "2.811.68"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "60.346.99"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "2.383.00"
$ws.Range("E18").Value = "  -2.83%  "
$ws.Range("D19").Value = "'10.61"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").Value = "'4.11"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("D21").Value = "'317.67"
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("D22").Value = "'6.69"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'1.81"
$ws.Range("E24").Value = "  +4.43%  "
$ws.Range("D25").Value = "'63.19"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "2.498.58"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").Value = "'7.80"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").Value = "0.0₃0926"
$ws.Range("E29").Value = "  -4.09%  "
$ws.Range("D30").Value = "'522.91"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  -3.75%  "
$ws.Range("D32").Value = "'8.00"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").Value = "'0.145"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'5.47"
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("D38").Value = "'4.65"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").Value = "'0.376"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'18.08"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").Value = "'1.73"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'138.14"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("D44").Value = "'40.18"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = "  -4.90%  "
$ws.Range("D46").Value = "'140.14"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "'3.54"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'20.38"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "'0.0515"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "'0.577"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'0.0927"
$ws.Range("E51").Value = "  -1.22%  "
